$wb = $excel.ActiveWorkbook

# Clear the "dt" interval column (D1:D7) on the EDTSlot sheet - the routine
# interval now lives in the routine config instead of per-slot.
$wsEdt = $wb.Worksheets.Item("EDTSlot")
$wsEdt.Activate()
$wsEdt.Range("D:D").Select()
$wsEdt.Range("D1:D7").ClearContents()

# Same cleanup on the UCTSlot sheet.
$wsUct = $wb.Worksheets.Item("UCTSlot")
$wsUct.Range("D1:D7").ClearContents()
$wsUct.Range("E28").Select()

# The rtnconfig model (RTEDCFG sheet) is no longer used - remove it entirely.
$wsCfg = $wb.Worksheets.Item("RTEDCFG")
$wsCfg.Delete()

$wsNsr = $wb.Worksheets.Item("NSR")
$wsNsr.Activate()
$wsNsr.Range("J19").Select()
